$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price/Volume columns (data is textual, not numeric)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "256.66"
$ws.Range("E2").Value = "-0.79%"
$ws.Range("D3").Value = "27.01"
$ws.Range("E3").Value = "-3.52%"
$ws.Range("D4").Value = "4.614"
$ws.Range("E4").Value = "-11.37%"
$ws.Range("D5").Value = "0.05906"
$ws.Range("E5").Value = "-0.07%"
$ws.Range("D6").Value = "6.624"
$ws.Range("D7").Value = "0.8601"
$ws.Range("E7").Value = "-1.64%"
$ws.Range("D8").Value = "0.9342"
$ws.Range("E8").Value = "-5.93%"
$ws.Range("E9").Value = "-0.24%"
$ws.Range("D10").Value = "0.03617"
$ws.Range("E10").Value = "-2.40%"
$ws.Range("D11").Value = "0.07089"
$ws.Range("E11").Value = "-1.36%"
$ws.Range("D12").Value = "0.03232"
$ws.Range("E12").Value = "2.22%"
$ws.Range("D13").Value = "0.09207"
$ws.Range("E13").Value = "-0.26%"
$ws.Range("D14").Value = "0.001548"
$ws.Range("E14").Value = "-0.08%"
$ws.Range("D15").Value = "0.0006071"
$ws.Range("E15").Value = "-94.30%"
$ws.Range("D16").Value = "0.006078"
$ws.Range("E16").Value = "0.97%"
$ws.Range("D17").Value = "3.516"
$ws.Range("E17").Value = "0.61%"
$ws.Range("D18").Value = "3.194"
$ws.Range("E18").Value = "-1.10%"
$ws.Range("E19").Value = "-0.13%"
$ws.Range("D20").Value = "0.3052"
$ws.Range("E20").Value = "-2.23%"
$ws.Range("D22").Value = "3.849"
$ws.Range("E22").Value = "9.10%"
$ws.Range("E23").Value = "1.00%"
$ws.Range("D24").Value = "0.001221"
$ws.Range("E24").Value = "0.20%"
$ws.Range("E25").Value = "-5.92%"
$ws.Range("E26").Value = "0.12%"
$ws.Range("E27").Value = "0.07%"
$ws.Range("D40").Value = "0.03829"
$ws.Range("E40").Value = "-0.26%"
$ws.Range("D41").Value = "0.006220"
$ws.Range("E41").Value = "14.06%"
$ws.Range("D42").Value = "0.1100"
$ws.Range("E42").Value = "-0.45%"
$ws.Range("D43").Value = "0.002200"
$ws.Range("E43").Value = "-4.23%"
$ws.Range("D44").Value = "0.01138"
$ws.Range("E44").Value = "6.83%"
$ws.Range("D45").Value = "0.00005457"
$ws.Range("E45").Value = "0.82%"
$ws.Range("E46").Value = "0.12%"
$ws.Range("E47").Value = "-29.50%"
$ws.Range("D48").Value = "0.09096"
$ws.Range("E48").Value = "4,160.08%"
$ws.Range("E49").Value = "0.12%"
$ws.Range("E50").Value = "0.12%"
